# Update Excel file with latest predictions
# The "EV Away win" sheet had its first data row (AUSTRALIA / Moreton City
# Excelsior - Peninsula Power) removed; the remaining prediction rows shift
# up by one, so the sheet's used range shrinks from A1:G6 to A1:G5.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("EV Away win")

# Remove the obsolete prediction row (row 2) — Excel shifts the rows below
# it up automatically, just like deleting a row via the UI / Range.Delete.
$ws.Rows.Item(2).Delete()
